$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")
$ws.Activate()

# New "Grid" PSET row appended to VEDA_Sets-Proc (PSET_SET / PSET_PN / SetName)
$ws.Range("F21").Value = "Grid"
$ws.Range("A21").Value = "IRE"
$ws.Range("B21").Value = "g[_]*"

$ws.Range("B21").Select()
